# Atualização de bases das ligas, do dia: 16-06-2024 às 07:16
#
# The source feed re-sorted a handful of fixture rows (same Date/Div block,
# adjacent match ids). This swaps the full record (every column except the
# running "id" in column A) between each pair of rows so the sheet reflects
# the corrected ordering coming from the upstream scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is a pair of row numbers whose B:AD contents (id in column A is
# left untouched) need to trade places.
$rowPairs = @(
    @(19, 20),
    @(90, 91),
    @(167, 168),
    @(258, 259),
    @(279, 280),
    @(291, 294)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
